{"js": "// The document starts with a Title paragraph, an Author paragraph, an\n// AbstractTitle paragraph, and an Abstract paragraph (in that order).\n// Each of the Title/Author/Abstract paragraphs originally had its text\n// split across many runs (one run per word/space). The edit simply\n// collapses each of those paragraphs down to a single run containing\n// the same overall text - no wording changes, just fewer runs.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Map each paragraph's current (multi-run) text to the merged text it\n// should become. Matching on text (trimmed) rather than a hard-coded\n// index keeps this robust to the exact paragraph position.\nconst replacements = [\n  {\n    from: \"Answers: Introduction to rearranging equations\",\n    to: \"Answers: Introduction to rearranging equations\",\n  },\n  {\n    from: \"Shanelle Advani, Tom Coleman\",\n    to: \"Shanelle Advani, Tom Coleman\",\n  },\n  {\n    from:\n      \"Answers to questions relating to the guide on introduction to rearranging equations.\",\n    to:\n      \"Answers to questions relating to the guide on introduction to rearranging equations.\",\n  },\n];\n\nfor (const paragraph of paragraphs.items) {\n  const current = paragraph.text.trim();\n  const match = replacements.find((r) => r.from === current);\n  if (match) {\n    // insertText(..., \"Replace\") rewrites the paragraph's range with a\n    // single new run, which is exactly the word-run-merge the diff shows.\n    paragraph.insertText(match.to, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# The document starts with a Title paragraph, an Author paragraph, an\n# AbstractTitle paragraph, and an Abstract paragraph (in that order).\n# Each of the Title/Author/Abstract paragraphs originally had its text\n# split across many runs (one run per word/space). The edit simply\n# collapses each of those paragraphs down to a single run containing\n# the same overall text - no wording changes.\n#\n# Find.Execute with Replace:=wdReplaceAll rewrites the whole matched\n# range as a single run, which is exactly the run-merging behaviour we\n# need (a plain `Range.Text = ...` assignment would only overwrite the\n# first run and leave the remaining word/space runs dangling).\n\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$titleText = \"Answers: Introduction to rearranging equations\"\n$titleRange = $d.Paragraphs(1).Range\n$null = $titleRange.Find.Execute($titleText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $titleText, $wdReplaceAll)\n\n$authorText = \"Shanelle Advani, Tom Coleman\"\n$authorRange = $d.Paragraphs(2).Range\n$null = $authorRange.Find.Execute($authorText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $authorText, $wdReplaceAll)\n\n$abstractText = \"Answers to questions relating to the guide on introduction to rearranging equations.\"\n$abstractRange = $d.Paragraphs(4).Range\n$null = $abstractRange.Find.Execute($abstractText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $abstractText, $wdReplaceAll)\n"}
